$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "77÷5=15, 2"
$tbl.Cell(1, 2).Range.Text = "35÷8=4, 3"
$tbl.Cell(1, 3).Range.Text = "59÷2=29, 1"
$tbl.Cell(1, 4).Range.Text = "40÷9=4, 4"
$tbl.Cell(1, 5).Range.Text = "96÷8=12, 0"

$tbl.Cell(5, 1).Range.Text = "77÷4=19, 1"
$tbl.Cell(5, 2).Range.Text = "40÷2=20, 0"
$tbl.Cell(5, 3).Range.Text = "80÷3=26, 2"
$tbl.Cell(5, 4).Range.Text = "71÷2=35, 1"
$tbl.Cell(5, 5).Range.Text = "25÷5=5, 0"

$tbl.Cell(9, 1).Range.Text = "47÷9=5, 2"
$tbl.Cell(9, 2).Range.Text = "63÷2=31, 1"
$tbl.Cell(9, 3).Range.Text = "58÷9=6, 4"
$tbl.Cell(9, 4).Range.Text = "41÷9=4, 5"
$tbl.Cell(9, 5).Range.Text = "63÷6=10, 3"

$tbl.Cell(13, 1).Range.Text = "31÷3=10, 1"
$tbl.Cell(13, 2).Range.Text = "52÷9=5, 7"
$tbl.Cell(13, 3).Range.Text = "37÷4=9, 1"
$tbl.Cell(13, 4).Range.Text = "61÷9=6, 7"
$tbl.Cell(13, 5).Range.Text = "88÷8=11, 0"

$tbl.Cell(17, 1).Range.Text = "30÷6=5, 0"
$tbl.Cell(17, 2).Range.Text = "16÷9=1, 7"
$tbl.Cell(17, 3).Range.Text = "33÷4=8, 1"
$tbl.Cell(17, 4).Range.Text = "90÷4=22, 2"
$tbl.Cell(17, 5).Range.Text = "46÷4=11, 2"
